# The "numberOfPages" column (column P) is being removed from Sheet1.
# Deleting the entire column shifts everything to its right one column to
# the left (Q->P, R->Q, ..., AE->AD), fixes up the dimension, the
# dataValidation sqref, and drops the now-unused "numberOfPages" shared
# string automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns("P").Delete()

# Leave the view focused on the column that now occupies the old P spot,
# matching the post-edit selection state.
$ws.Columns("P").Select()
